# code tidy up and new examples.
# Insert a new blank row above row 2 on the "Tree 1" and "Tree 2" sheets,
# pushing all existing bracket content down by one row.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Tree 1")
$ws1.Rows("2:2").Insert()

$ws2 = $wb.Worksheets.Item("Tree 2")
$ws2.Rows("2:2").Insert()
